$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted above the existing data
# (row 77), pushing all subsequent rows (old 77-137) down by one to
# (78-138). Insert a fresh row 77 so the rest of the table shifts down
# intact, then populate it with the new record's values.
$ws.Rows("77:77").Insert()

$ws.Cells.Item(77, 1).Value = 7
$ws.Cells.Item(77, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(77, 3).Value = "Ñuble"
$ws.Cells.Item(77, 4).Value = 45072
$ws.Cells.Item(77, 5).Value = 16
$ws.Cells.Item(77, 6).Value = 100112031
$ws.Cells.Item(77, 7).Value = "Poroto verde"
$ws.Cells.Item(77, 8).Value = "Magnum"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 20
$ws.Cells.Item(77, 11).Value = 35000
$ws.Cells.Item(77, 12).Value = 35000
$ws.Cells.Item(77, 13).Value = 35000
$ws.Cells.Item(77, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(77, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(77, 16).Value = 1400
$ws.Cells.Item(77, 17).Value = 25
$ws.Cells.Item(77, 18).Value = "Hortaliza"
